$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = 2.5
    $cell.Borders.LineStyle = -4142
}

$ws.Range("H9").Select()
